$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Sample Number" values for rows 82-97 (A82:A97 = 81..96)
for ($r = 82; $r -le 97; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Paste the CSV library-input sample names into column D (D2:D97 = AAA01..AAA96)
for ($r = 2; $r -le 97; $r++) {
    $n = $r - 1
    $ws.Cells.Item($r, 4).Value = "AAA" + $n.ToString("00")
}

$ws.Columns("D").EntireColumn.AutoFit()

$ws.Range("C96").Select()
